# Fixed hidden slides on presentation.
# Slides 2, 3, and 8 were marked as hidden (show="0" in the slide XML).
# Un-hide them by clearing the SlideShowTransition.Hidden flag.

$p = $ppt.ActivePresentation

foreach ($idx in 2, 3, 8) {
    $s = $p.Slides.Item($idx)
    $s.SlideShowTransition.Hidden = 0
}
